$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("RUIZ RUIZ LUZ MERI", 69),
    @("VASQUEZ DIAZ LUZ ANGELICA", 68),
    @("BENAVIDES MARRUFO ARACELI", 66),
    @("SOTO VALLEJOS ELSITA", 66),
    @("ROJAS VASQUEZ FLOR NOELITA", 64),
    @("MEDINA VALLEJOS ERICK LEONARDO", 64),
    @("ZAMORA TAMAY NEYSER IVAN", 63),
    @("TELLO FERNANDEZ MILENY", 63),
    @("TIRADO PEREZ JEINER", 62),
    @("SOTO VILLENA NILSON", 62),
    @("PÓSITO CHUGDEN NANIX", 61),
    @("GALLARDO CORTEZ MELISSA DEL CARMEN", 56),
    @("VASQUEZ LUNA YUDITH", 56),
    @("BENAVIDES SALAZAR IDELSA", 49)
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $row = $row + 1
}
